$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')

# Insert the new row (shifts everything from row 17 down by one)
$ws.Rows.Item(17).Insert()

# Match column-A formatting (bold centered text with thin border) used throughout the table
$styleSrc = $ws.Range("A16")
$styleDst = $ws.Range("A17")
$styleDst.Font.Bold = $styleSrc.Font.Bold
$styleDst.HorizontalAlignment = $styleSrc.HorizontalAlignment
$styleDst.VerticalAlignment = $styleSrc.VerticalAlignment
$styleDst.Borders.LineStyle = $styleSrc.Borders.LineStyle

# Populate every data row in this block with the target (post-edit) content
$ws.Range("A17").Value = 16
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '2024-06-22'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '萍乡·AU9夏至国漫展'
$ws.Range("D17").Value = '金陵东路18号 萍乡市体育馆'
$ws.Range("E17").Value = '2024.06.22 10:00-06.22 17:00'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=86453'
$ws.Range("I17").Value = '//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg'

$ws.Range("A18").Value = 17
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2024-06-23'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '上饶·BM次元盛典运动番only'
$ws.Range("D18").Value = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws.Range("E18").Value = '2024.06.23 10:00-06.23 17:00'
$ws.Range("F18").Value = 215
$ws.Range("G18").Value = 55
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=85201'
$ws.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png'

$ws.Range("A19").Value = 18
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '2024-06-29'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '萍乡·BM次元盛典运动番only'
$ws.Range("D19").Value = '康庄路3号 萍乡梅园国际大酒店'
$ws.Range("E19").Value = '2024.06.29 10:00-06.29 17:00'
$ws.Range("F19").Value = 222
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=85192'
$ws.Range("I19").Value = '//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

$ws.Range("A20").Value = 19
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '2024-06-30'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '宜春·BM次元盛典运动番only'
$ws.Range("D20").Value = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws.Range("E20").Value = '2024.06.30 10:00-06.30 17:00'
$ws.Range("F20").Value = 216
$ws.Range("G20").Value = 55
$ws.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=84636'
$ws.Range("I20").Value = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

$ws.Range("A21").Value = 20
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '2024-07-06'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '鹰潭·BM次元盛典运动番only'
$ws.Range("D21").Value = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws.Range("E21").Value = '2024.07.06 10:00-07.06 17:00'
$ws.Range("F21").Value = 13
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=85997'
$ws.Range("I21").Value = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

$ws.Range("A22").Value = 21
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '2024-07-14'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '吉安·COMIC LIFE次元假日05'
$ws.Range("D22").Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Range("E22").Value = '2024.07.14 09:00-07.14 18:00'
$ws.Range("F22").Value = 198
$ws.Range("G22").Value = 9.9
$ws.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=85924'
$ws.Range("I22").Value = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

$ws.Range("A23").Value = 22
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2024-07-20'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '南昌·漫拥动漫嘉年华Pro-追光启航'
$ws.Range("D23").Value = '小蓝南路420号 洪州体育馆'
$ws.Range("E23").Value = '2024.07.20 09:00-07.21 17:00'
$ws.Range("F23").Value = 65
$ws.Range("G23").Value = 52.5
$ws.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=85796'
$ws.Range("I23").Value = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

$ws.Range("A24").Value = 23
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2024-07-27'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '江西·次元星河国风动漫游戏嘉年华'
$ws.Range("D24").Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Range("E24").Value = '2024.07.27 10:00-07.28 17:00'
$ws.Range("F24").Value = 1709
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=85493'
$ws.Range("I24").Value = '//i2.hdslb.com/bfs/openplatform/202404/HJ7TF5zx1714367786872.jpeg'

$ws.Range("A25").Value = 24
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2024-07-28'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws.Range("D25").Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$ws.Range("E25").Value = '2024.07.28 11:00-07.28 17:00'
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=85688'
$ws.Range("I25").Value = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

$ws.Range("A26").Value = 25
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '2024-08-03'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '南昌·幻梦境国际动漫游戏嘉年华1th'
$ws.Range("D26").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("E26").Value = '2024.08.03 09:00-08.04 17:30'
$ws.Range("F26").Value = 416
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=83980'
$ws.Range("I26").Value = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

$ws.Range("A27").Value = 26
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '2024-08-03'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws.Range("D27").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("E27").Value = '2024.08.03 09:00-08.04 17:00'
$ws.Range("F27").Value = 43
$ws.Range("G27").Value = 55
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=86341'
$ws.Range("I27").Value = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

$ws.Range("A28").Value = 27
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '2024-08-03'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws.Range("D28").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("E28").Value = '2024.08.03 08:30-08.03 17:00'
$ws.Range("F28").Value = 571
$ws.Range("G28").Value = '已售罄'
$ws.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=85981'
$ws.Range("I28").Value = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

$ws.Range("A29").Value = 28
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '2024-08-04'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '九江·第一届异次元动漫嘉年华'
$ws.Range("D29").Value = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws.Range("E29").Value = '2024.08.04 08:00-08.04 17:00'
$ws.Range("F29").Value = 213
$ws.Range("G29").Value = 45
$ws.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=84407'
$ws.Range("I29").Value = '//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg'

$ws.Range("A30").Value = 29
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '2024-08-06'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '南昌·第一届异次元动漫嘉年华'
$ws.Range("D30").Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws.Range("E30").Value = '2024.08.06 08:00-08.06 17:00'
$ws.Range("F30").Value = 306
$ws.Range("G30").Value = 55
$ws.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=84102'
$ws.Range("I30").Value = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

$ws.Range("A31").Value = 30
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = '2024-08-08'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '赣州·第二届异次元动漫嘉年华'
$ws.Range("D31").Value = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws.Range("E31").Value = '2024.08.08 08:00-08.08 17:00'
$ws.Range("F31").Value = 434
$ws.Range("G31").Value = 45
$ws.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=84184'
$ws.Range("I31").Value = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

# "Want to go" counter bumps for rows outside the inserted block
$ws.Range("F2").Value = 1893
$ws.Range("F3").Value = 503
$ws.Range("F6").Value = 2665
$ws.Range("F7").Value = 182
$ws.Range("F10").Value = 1560
$ws.Range("F11").Value = 544

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')

# Insert the new row (shifts everything from row 18 down by one)
$ws.Rows.Item(18).Insert()

# Match column-A formatting (bold centered text with thin border) used throughout the table
$styleSrc = $ws.Range("A17")
$styleDst = $ws.Range("A18")
$styleDst.Font.Bold = $styleSrc.Font.Bold
$styleDst.HorizontalAlignment = $styleSrc.HorizontalAlignment
$styleDst.VerticalAlignment = $styleSrc.VerticalAlignment
$styleDst.Borders.LineStyle = $styleSrc.Borders.LineStyle

# Populate every data row in this block with the target (post-edit) content
$ws.Range("A18").Value = 16
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2024-06-22'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '萍乡·AU9夏至国漫展'
$ws.Range("D18").Value = '金陵东路18号 萍乡市体育馆'
$ws.Range("E18").Value = '2024.06.22 10:00-06.22 17:00'
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 45
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=86453'
$ws.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg'

$ws.Range("A19").Value = 17
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '2024-06-23'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '上饶·BM次元盛典运动番only'
$ws.Range("D19").Value = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws.Range("E19").Value = '2024.06.23 10:00-06.23 17:00'
$ws.Range("F19").Value = 215
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=85201'
$ws.Range("I19").Value = '//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png'

$ws.Range("A20").Value = 18
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '2024-06-29'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '萍乡·BM次元盛典运动番only'
$ws.Range("D20").Value = '康庄路3号 萍乡梅园国际大酒店'
$ws.Range("E20").Value = '2024.06.29 10:00-06.29 17:00'
$ws.Range("F20").Value = 222
$ws.Range("G20").Value = 55
$ws.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=85192'
$ws.Range("I20").Value = '//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

$ws.Range("A21").Value = 19
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '2024-06-30'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '宜春·BM次元盛典运动番only'
$ws.Range("D21").Value = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws.Range("E21").Value = '2024.06.30 10:00-06.30 17:00'
$ws.Range("F21").Value = 216
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=84636'
$ws.Range("I21").Value = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

$ws.Range("A22").Value = 20
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '2024-07-06'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '鹰潭·BM次元盛典运动番only'
$ws.Range("D22").Value = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws.Range("E22").Value = '2024.07.06 10:00-07.06 17:00'
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = 55
$ws.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=85997'
$ws.Range("I22").Value = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

$ws.Range("A23").Value = 21
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2024-07-14'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '吉安·COMIC LIFE次元假日05'
$ws.Range("D23").Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Range("E23").Value = '2024.07.14 09:00-07.14 18:00'
$ws.Range("F23").Value = 198
$ws.Range("G23").Value = 9.9
$ws.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=85924'
$ws.Range("I23").Value = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

$ws.Range("A24").Value = 22
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2024-07-20'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '南昌·漫拥动漫嘉年华Pro-追光启航'
$ws.Range("D24").Value = '小蓝南路420号 洪州体育馆'
$ws.Range("E24").Value = '2024.07.20 09:00-07.21 17:00'
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 52.5
$ws.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=85796'
$ws.Range("I24").Value = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

$ws.Range("A25").Value = 23
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2024-07-27'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '江西·次元星河国风动漫游戏嘉年华'
$ws.Range("D25").Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Range("E25").Value = '2024.07.27 10:00-07.28 17:00'
$ws.Range("F25").Value = 1709
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=85493'
$ws.Range("I25").Value = '//i2.hdslb.com/bfs/openplatform/202404/HJ7TF5zx1714367786872.jpeg'

$ws.Range("A26").Value = 24
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '2024-07-28'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws.Range("D26").Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$ws.Range("E26").Value = '2024.07.28 11:00-07.28 17:00'
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=85688'
$ws.Range("I26").Value = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

$ws.Range("A27").Value = 25
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '2024-08-03'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = '南昌·幻梦境国际动漫游戏嘉年华1th'
$ws.Range("D27").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("E27").Value = '2024.08.03 09:00-08.04 17:30'
$ws.Range("F27").Value = 416
$ws.Range("G27").Value = 64
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=83980'
$ws.Range("I27").Value = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

$ws.Range("A28").Value = 26
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '2024-08-03'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws.Range("D28").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("E28").Value = '2024.08.03 09:00-08.04 17:00'
$ws.Range("F28").Value = 43
$ws.Range("G28").Value = 55
$ws.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=86341'
$ws.Range("I28").Value = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

$ws.Range("A29").Value = 27
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '2024-08-03'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws.Range("D29").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("E29").Value = '2024.08.03 08:30-08.03 17:00'
$ws.Range("F29").Value = 571
$ws.Range("G29").Value = '已售罄'
$ws.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=85981'
$ws.Range("I29").Value = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

$ws.Range("A30").Value = 28
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '2024-08-04'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '九江·第一届异次元动漫嘉年华'
$ws.Range("D30").Value = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws.Range("E30").Value = '2024.08.04 08:00-08.04 17:00'
$ws.Range("F30").Value = 213
$ws.Range("G30").Value = 45
$ws.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=84407'
$ws.Range("I30").Value = '//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg'

$ws.Range("A31").Value = 29
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = '2024-08-06'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '南昌·第一届异次元动漫嘉年华'
$ws.Range("D31").Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws.Range("E31").Value = '2024.08.06 08:00-08.06 17:00'
$ws.Range("F31").Value = 306
$ws.Range("G31").Value = 55
$ws.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=84102'
$ws.Range("I31").Value = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

$ws.Range("A32").Value = 30
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = '2024-08-08'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = '赣州·第二届异次元动漫嘉年华'
$ws.Range("D32").Value = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws.Range("E32").Value = '2024.08.08 08:00-08.08 17:00'
$ws.Range("F32").Value = 434
$ws.Range("G32").Value = 45
$ws.Range("H32").Value = 'https://show.bilibili.com/platform/detail.html?id=84184'
$ws.Range("I32").Value = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

# "Want to go" counter bumps for rows outside the inserted block
$ws.Range("F2").Value = 1893
$ws.Range("F4").Value = 503
$ws.Range("F7").Value = 2665
$ws.Range("F8").Value = 182
$ws.Range("F11").Value = 1560
$ws.Range("F12").Value = 544

